$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 41, shifting existing rows 41..186 down to 42..187
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new data record
$ws.Cells.Item(41, 1).Value = 8
$ws.Cells.Item(41, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(41, 3).Value = "Coquimbo"
$ws.Cells.Item(41, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
$ws.Cells.Item(41, 4).Value = 44453
$ws.Cells.Item(41, 5).Value = 4
$ws.Cells.Item(41, 6).Value = 100114013
$ws.Cells.Item(41, 7).Value = "Zanahoria"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 700
$ws.Cells.Item(41, 11).Value = 4500
$ws.Cells.Item(41, 12).Value = 5000
$ws.Cells.Item(41, 13).Value = 4750
$ws.Cells.Item(41, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(41, 16).Value = 238
$ws.Cells.Item(41, 17).Value = 20
$ws.Cells.Item(41, 18).Value = "Hortaliza"
